# This script updates the "江西-漫展信息" workbook to reflect the new
# scrape snapshot:
#  - The 宜春 event (previously row 2) has been removed entirely from
#    both the "展览" sheet and the combined "全部类型" sheet, shifting
#    every subsequent row up by one.
#  - Various "想去人数" (F column) interest counters were refreshed to
#    their latest values (including the 南昌·Kpop New Life row, which
#    lives on its own on the "演出" sheet as well as inside "全部类型").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览": delete the 宜春 row (row 2), then refresh F values
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows.Item(2).Delete()

$ws1.Range("F2").Value = 873
$ws1.Range("F4").Value = 256
$ws1.Range("F6").Value = 161
$ws1.Range("F7").Value = 147
$ws1.Range("F9").Value = 4555
$ws1.Range("F12").Value = 500
$ws1.Range("F13").Value = 456
$ws1.Range("F14").Value = 17
$ws1.Range("F16").Value = 1234
$ws1.Range("F17").Value = 2495
$ws1.Range("F18").Value = 377
$ws1.Range("F19").Value = 73
$ws1.Range("F20").Value = 52
$ws1.Range("F22").Value = 2264
$ws1.Range("F23").Value = 91
$ws1.Range("F25").Value = 26
$ws1.Range("F26").Value = 159
$ws1.Range("F29").Value = 223
$ws1.Range("F30").Value = 35

# ---------------------------------------------------------------
# Sheet "演出": refresh the 南昌·Kpop New Life interest counter
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 38

# ---------------------------------------------------------------
# Sheet "全部类型": delete the 宜春 row (row 2), then refresh F values
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(2).Delete()

$ws4.Range("F2").Value = 873
$ws4.Range("F4").Value = 256
$ws4.Range("F6").Value = 161
$ws4.Range("F7").Value = 147
$ws4.Range("F9").Value = 38
$ws4.Range("F10").Value = 4555
$ws4.Range("F13").Value = 500
$ws4.Range("F14").Value = 456
$ws4.Range("F15").Value = 17
$ws4.Range("F17").Value = 1234
$ws4.Range("F18").Value = 2495
$ws4.Range("F19").Value = 377
$ws4.Range("F20").Value = 73
$ws4.Range("F21").Value = 52
$ws4.Range("F23").Value = 2264
$ws4.Range("F24").Value = 91
$ws4.Range("F26").Value = 26
$ws4.Range("F27").Value = 159
$ws4.Range("F30").Value = 223
$ws4.Range("F31").Value = 35

$wb.Save()
